$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the program_name column (D2:D4) from "Toyota Rewards Visa" to "Auto Rewards Visa"
$ws.Range("D2:D4").Value = "Auto Rewards Visa"

# Reflect the active selection recorded after making this edit
$ws.Range("D2:D4").Select()
